$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "contoh"
$ws.Range("A3:B5").ClearContents()
$ws.Range("E3").Select()
